# Insert a new data row at row 7 (pushes existing rows 7..81 down to 8..82)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value  = 1
$ws.Cells.Item(7, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value  = 44545
$ws.Cells.Item(7, 5).Value  = 15
$ws.Cells.Item(7, 6).Value  = "Fruta"
$ws.Cells.Item(7, 7).Value  = 100102
$ws.Cells.Item(7, 8).Value  = "Cítricos"
$ws.Cells.Item(7, 9).Value  = 100102004
$ws.Cells.Item(7, 10).Value = "Mandarina"
$ws.Cells.Item(7, 11).Value = "Murcott"
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 300
$ws.Cells.Item(7, 14).Value = 13000
$ws.Cells.Item(7, 15).Value = 14000
$ws.Cells.Item(7, 16).Value = 13500
$ws.Cells.Item(7, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(7, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(7, 19).Value = 675
$ws.Cells.Item(7, 20).Value = 20
